$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.192.11"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "3.070.58"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'558.11"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "'145.78"
$ws.Range("E6").Value = "  +4.93%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.064.77"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.152"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'6.27"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("D13").Value = "'0.0000227"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'35.12"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "3.583.43"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "64.218.92"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "3.077.18"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "'6.77"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "'475.28"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'13.93"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").Value = "'0.675"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'7.53"
$ws.Range("E23").Value = "  +4.34%  "
$ws.Range("D24").Value = "'13.54"
$ws.Range("E24").Value = "  +7.80%  "
$ws.Range("D25").Value = "'81.39"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'2.79"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").Value = "'8.09"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "'2.05"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'26.13"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "'2.48"
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("D34").Value = "'5.57"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").Value = "'6.16"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").Value = "'54.65"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "'463.08"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "'3.00"
$ws.Range("E38").Value = "  +14.26%  "
$ws.Range("D39").Value = "'0.0830"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("D40").Value = "'0.0404"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("D41").Value = "2.963.18"
$ws.Range("E41").Value = "  -5.99%  "
$ws.Range("D42").Value = "'8.26"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'0.114"
$ws.Range("E43").Value = "  -4.62%  "
$ws.Range("D44").Value = "'28.03"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "'0.259"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.13"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").Value = "'0.112"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").Value = "'119.66"
$ws.Range("E49").Value = "  +3.30%  "
$ws.Range("D50").Value = "0.0₃0516"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'2.06"
$ws.Range("E51").Value = "  +0.03%  "
